# The workbook has two sheets: "2020-Q4" (first tab) and "总计" (second tab).
# Re-sort the sheet tabs so that "总计" (the summary sheet) comes first,
# followed by "2020-Q4" (the detail sheet) - i.e. swap their order/positions.
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)

$summarySheet.Move($firstSheet)

# Make sure the newly-first sheet ("总计") is the active/selected tab,
# matching the workbook's saved active-tab state.
$wb.Worksheets.Item("总计").Activate()
